$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (D) and Volume(1h) (E) values from the latest scrape.
$values = @{
    2 = @{ "D" = "29.381.34"; "E" = "  -0.25%  " }
    3 = @{ "D" = "1.846.87"; "E" = "  -0.32%  " }
    4 = @{ "D" = "1.000"; "E" = "  +0.07%  " }
    5 = @{ "D" = "238.76"; "E" = "  -1.47%  " }
    6 = @{ "D" = "0.6299"; "E" = "  -0.17%  " }
    7 = @{ "E" = "  +0.08%  " }
    8 = @{ "D" = "0.07543"; "E" = "  -0.67%  " }
    9 = @{ "D" = "0.2946"; "E" = "  -1.19%  " }
    10 = @{ "D" = "24.65"; "E" = "  +0.38%  " }
    11 = @{ "D" = "0.07699"; "E" = "  -0.10%  " }
    12 = @{ "D" = "1.864.00"; "E" = "  -3.09%  " }
    13 = @{ "D" = "4.987"; "E" = "  -0.36%  " }
    14 = @{ "D" = "0.6798"; "E" = "  -1.45%  " }
    15 = @{ "D" = "0.00001017"; "E" = "  +2.52%  " }
    16 = @{ "D" = "83.16"; "E" = "  -0.38%  " }
    17 = @{ "D" = "2.109.95"; "E" = "  -3.82%  " }
    18 = @{ "D" = "6.133"; "E" = "  -0.72%  " }
    19 = @{ "D" = "29.429.80"; "E" = "  -0.59%  " }
    20 = @{ "D" = "227.84"; "E" = "  -2.35%  " }
    21 = @{ "D" = "12.48"; "E" = "  -0.38%  " }
    22 = @{ "E" = "  -0.02%  " }
    23 = @{ "D" = "7.501"; "E" = "  -2.68%  " }
    24 = @{ "D" = "1.002"; "E" = "  +0.18%  " }
    25 = @{ "D" = "156.90"; "E" = "  +1.13%  " }
    26 = @{ "D" = "0.1395"; "E" = "  -0.16%  " }
    27 = @{ "D" = "8.369"; "E" = "  -1.23%  " }
    28 = @{ "D" = "17.62"; "E" = "  -0.41%  " }
    29 = @{ "D" = "1.463"; "E" = "  -0.89%  " }
    30 = @{ "D" = "1.273"; "E" = "  +0.67%  " }
    31 = @{ "D" = "0.05667"; "E" = "  -2.30%  " }
    32 = @{ "D" = "4.124"; "E" = "  -0.12%  " }
    33 = @{ "D" = "4.031"; "E" = "  +0.33%  " }
    34 = @{ "D" = "1.838"; "E" = "  -2.92%  " }
    35 = @{ "D" = "1.156"; "E" = "  -1.19%  " }
    36 = @{ "D" = "0.7152"; "E" = "  -1.28%  " }
    37 = @{ "E" = "  +0.12%  " }
    38 = @{ "D" = "1.247.89"; "E" = "  -0.70%  " }
    39 = @{ "D" = "0.01805"; "E" = "  -0.02%  " }
    40 = @{ "D" = "2.775"; "E" = "  -0.80%  " }
    41 = @{ "D" = "6.191"; "E" = "  +1.20%  " }
    42 = @{ "D" = "0.9047"; "E" = "  -0.69%  " }
    43 = @{ "E" = "  +0.11%  " }
    44 = @{ "D" = "101.58"; "E" = "  -0.16%  " }
    45 = @{ "D" = "66.26"; "E" = "  -2.25%  " }
    46 = @{ "D" = "0.00000000118"; "E" = "  +1.13%  " }
    47 = @{ "D" = "7.068"; "E" = "  -3.62%  " }
    48 = @{ "D" = "0.4014"; "E" = "  -0.64%  " }
    49 = @{ "D" = "9.048"; "E" = "  -1.60%  " }
    50 = @{ "D" = "1.691"; "E" = "  -1.06%  " }
    51 = @{ "D" = "0.1116"; "E" = "  -0.21%  " }
}

# Rows whose Price text looks like a plain number (e.g. "1.000") - these must
# be pre-formatted as Text so Excel keeps them as literal strings instead of
# coercing/rounding them into numeric values.
$textFormatRows = @(4, 5, 6, 8, 9, 10, 11, 13, 14, 15, 16, 18, 20, 21, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 39, 40, 41, 42, 44, 45, 46, 47, 48, 49, 50, 51)
foreach ($row in $textFormatRows) {
    $ws.Cells.Item($row, 4).NumberFormat = "@"
}

foreach ($row in $values.Keys) {
    $cells = $values[$row]
    if ($cells.ContainsKey("D")) {
        $ws.Cells.Item($row, 4).Value = $cells["D"]
    }
    if ($cells.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = $cells["E"]
    }
}
